$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.227.96'
$ws.Range('E2').Value = '  +3.60%  '

$ws.Range('D3').Value = '3.072.64'
$ws.Range('E3').Value = '  +6.42%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').Value = "'514.92"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.80%  '

$ws.Range('D6').Value = "'140.68"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +7.41%  '

$ws.Range('E7').Value = '  +0.16%  '

$ws.Range('D8').Value = "'0.433"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.32%  '

$ws.Range('D9').Value = "'7.23"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.87%  '

$ws.Range('E10').Value = '  +5.49%  '

$ws.Range('E11').Value = '  +7.95%  '

$ws.Range('D12').Value = '3.590.59'
$ws.Range('E12').Value = '  +6.22%  '

$ws.Range('E13').Value = '  +3.22%  '

$ws.Range('D14').Value = "'25.39"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.48%  '

$ws.Range('E15').Value = '  +5.66%  '

$ws.Range('D16').Value = '57.261.55'
$ws.Range('E16').Value = '  +3.75%  '

$ws.Range('D17').Value = '3.070.08'
$ws.Range('E17').Value = '  +6.37%  '

$ws.Range('D18').Value = "'5.92"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.41%  '

$ws.Range('D19').Value = "'13.02"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.53%  '

$ws.Range('D20').Value = "'8.16"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +8.24%  '

$ws.Range('D21').Value = "'336.90"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +8.59%  '

$ws.Range('E22').Value = '  +0.25%  '

$ws.Range('E23').Value = '  +5.40%  '

$ws.Range('D24').Value = "'65.20"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.60%  '

$ws.Range('D25').Value = "'0.170"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +8.13%  '

$ws.Range('E26').Value = '  +14.96%  '

$ws.Range('E27').Value = '  +0.21%  '

$ws.Range('D28').Value = "'6.46"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.58%  '

$ws.Range('D29').Value = "'7.08"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.50%  '

$ws.Range('E30').Value = '  +4.85%  '

$ws.Range('D31').Value = "'20.72"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.55%  '

$ws.Range('E32').Value = '  +7.58%  '

$ws.Range('E33').Value = '  +4.32%  '

$ws.Range('E34').Value = '  +5.65%  '

$ws.Range('E35').Value = '  +6.40%  '

$ws.Range('D36').Value = "'26.28"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +9.36%  '

$ws.Range('E37').Value = '  +6.04%  '

$ws.Range('D38').Value = "'0.0671"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.57%  '

$ws.Range('D39').Value = '3.106.44'
$ws.Range('E39').Value = '  +6.44%  '

$ws.Range('D40').Value = "'36.96"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.45%  '

$ws.Range('D41').Value = "'0.669"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.74%  '

$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = "'3.84"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.10%  '

$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = "'1.00"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.05%  '

$ws.Range('D44').Value = '2.237.55'
$ws.Range('E44').Value = '  +8.09%  '

$ws.Range('E45').Value = '  +11.16%  '

$ws.Range('E46').Value = '  +5.09%  '

$ws.Range('D47').Value = "'0.947"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.17%  '

$ws.Range('D48').Value = "'19.94"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +8.92%  '

$ws.Range('E49').Value = '  +0.70%  '

$ws.Range('E50').Value = '  +4.15%  '

$ws.Range('B51').Value = 'TheGraph'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D51').Value = "'0.181"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.26%  '
